$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "1.3.1.1f Доля лиц, получающих пенсии и пособия по инвалидности к общей численности населения"

$ws.Range("S2:S5").Copy()
$ws.Range("T2:T5").PasteSpecial()

$ws.Range("T2").Value = ""
$ws.Range("T3").Value = 2023
$ws.Range("T4").Value = 217222
$ws.Range("T5").Value = 2.9794303052841493
